$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.881.83'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.442.88'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.09%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '575.30'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -1.24%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '147.06'
$cell.ClearFormats()
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").Value = '3.443.16'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.01%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.476'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +0.08%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '7.76'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("E11").Value = '  -1.56%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.405'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +2.83%  '
$ws.Range("D13").Value = '4.033.18'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("E14").Value = '  +2.59%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '28.68'
$cell.ClearFormats()
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").Value = '3.440.08'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '62.991.56'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E19").Value = '  +1.55%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '14.30'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +0.22%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '9.12'
$cell.ClearFormats()
$ws.Range("E21").Value = '  -2.44%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '384.10'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -2.66%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.557'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -0.42%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '74.28'
$cell.ClearFormats()
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '3.574.08'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -4.16%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.181'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -5.88%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '7.58'
$cell.ClearFormats()
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("E30").Value = '  +0.10%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.00'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -1.85%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '2.09'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("E33").Value = '  -0.10%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '23.20'
$cell.ClearFormats()
$ws.Range("E34").Value = '  -2.03%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.27'
$cell.ClearFormats()
$ws.Range("E35").Value = '  -9.44%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '5.30'
$cell.ClearFormats()
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.59'
$cell.ClearFormats()
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '7.02'
$cell.ClearFormats()
$ws.Range("E38").Value = '  -1.57%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '31.73'
$cell.ClearFormats()
$ws.Range("E39").Value = '  +1.54%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '167.65'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("D41").Value = '3.479.05'
$ws.Range("E41").Value = '  +0.05%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.0762'
$cell.ClearFormats()
$ws.Range("E42").Value = '  -0.76%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.790'
$cell.ClearFormats()
$ws.Range("E43").Value = '  -0.58%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '42.32'
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.72'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").Value = '2.562.09'
$ws.Range("E48").Value = '  -0.23%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '2.28'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +2.88%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '6.83'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +1.56%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '22.45'
$cell.ClearFormats()
$ws.Range("E51").Value = '  -4.54%  '
